# Commit: Wed, Apr 15, 2020  5:04:57 AM
#
# 1) Three tables (on the slides that still used the bespoke "Table_0"
#    table style) get re-styled onto the built-in PowerPoint table style
#    {B00C50C8-98B5-4515-8BBF-F2E6FFDB2B85}.
# 2) The deck's applied Design swaps its colour scheme from the
#    "Integral / Red Violet" palette to the stock "Office" palette
#    (picking the built-in "Office Theme" colours from the Design
#    gallery) - the font scheme / effect scheme are unaffected because
#    both themes already shared them.

$p = $ppt.ActivePresentation

# --- 1. Retarget the three tables onto the built-in table style -----------
$builtInTableStyle = "{B00C50C8-98B5-4515-8BBF-F2E6FFDB2B85}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($builtInTableStyle, $false)
        }
    }
}

# --- 2. Switch the Design's colour scheme to the "Office" palette ---------
$colorScheme = $p.SlideMaster.ColorScheme
$colorScheme.Colors(1).RGB  = 0          # dk1      000000
$colorScheme.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$colorScheme.Colors(3).RGB  = 6968388    # dk2      44546A
$colorScheme.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$colorScheme.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$colorScheme.Colors(6).RGB  = 3243501    # accent2  ED7D31
$colorScheme.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$colorScheme.Colors(8).RGB  = 49407      # accent4  FFC000
$colorScheme.Colors(9).RGB  = 12874308   # accent5  4472C4
$colorScheme.Colors(10).RGB = 4697456    # accent6  70AD47
$colorScheme.Colors(11).RGB = 12673797   # hlink    0563C1
$colorScheme.Colors(12).RGB = 7491477    # folHlink 954F72
